$d = $word.ActiveDocument

# The document contains three tables:
#   1) "Verantwortlich / Information / ..." - keep this one.
#   2) "Zeile 1, Spalte 1 / ..."             - remove.
#   3) "Kopfzeile, Spalte 1 / ..."           - remove.
# Between tables 2 and 3 there is a single blank paragraph that also needs
# to be removed, while the blank paragraphs right before table 2 and right
# after table 3 must be kept untouched.

# Step 1: delete the blank paragraph sitting between tables 2 and 3. Both
# table lookups are done fresh (not cached across the subsequent deletes)
# since this runtime's Table references can become stale once other tables
# in the document are later deleted.
$gapStart = $d.Tables.Item(2).Range.End
$gapEnd = $d.Tables.Item(3).Range.Start
$d.Range($gapStart, $gapEnd).Delete()

# Step 2: remove the two sample tables. Fetch each one fresh immediately
# before deleting it, and go in descending index order so that removing
# table 3 first doesn't shift the index of table 2.
$d.Tables.Item(3).Delete()
$d.Tables.Item(2).Delete()
